$d = $word.ActiveDocument

$d.Content.Find.Execute("481×7=3367", $true, $false, $false, $false, $false, $true, 1, $false, "272×9=2448", 2)
$d.Content.Find.Execute("982×9=8838", $true, $false, $false, $false, $false, $true, 1, $false, "199×2=398", 2)
$d.Content.Find.Execute("241×4=964", $true, $false, $false, $false, $false, $true, 1, $false, "165×3=495", 2)
$d.Content.Find.Execute("959×3=2877", $true, $false, $false, $false, $false, $true, 1, $false, "244×2=488", 2)
$d.Content.Find.Execute("428×9=3852", $true, $false, $false, $false, $false, $true, 1, $false, "732×3=2196", 2)
$d.Content.Find.Execute("368×8=2944", $true, $false, $false, $false, $false, $true, 1, $false, "636×8=5088", 2)
$d.Content.Find.Execute("548×7=3836", $true, $false, $false, $false, $false, $true, 1, $false, "356×9=3204", 2)
$d.Content.Find.Execute("738×6=4428", $true, $false, $false, $false, $false, $true, 1, $false, "566×7=3962", 2)
$d.Content.Find.Execute("644×9=5796", $true, $false, $false, $false, $false, $true, 1, $false, "326×9=2934", 2)
$d.Content.Find.Execute("505×5=2525", $true, $false, $false, $false, $false, $true, 1, $false, "149×3=447", 2)
$d.Content.Find.Execute("530×2=1060", $true, $false, $false, $false, $false, $true, 1, $false, "332×4=1328", 2)
$d.Content.Find.Execute("585×8=4680", $true, $false, $false, $false, $false, $true, 1, $false, "768×5=3840", 2)
$d.Content.Find.Execute("313×9=2817", $true, $false, $false, $false, $false, $true, 1, $false, "493×7=3451", 2)
$d.Content.Find.Execute("198×2=396", $true, $false, $false, $false, $false, $true, 1, $false, "114×9=1026", 2)
$d.Content.Find.Execute("989×6=5934", $true, $false, $false, $false, $false, $true, 1, $false, "804×4=3216", 2)
$d.Content.Find.Execute("240×3=720", $true, $false, $false, $false, $false, $true, 1, $false, "293×8=2344", 2)
$d.Content.Find.Execute("386×2=772", $true, $false, $false, $false, $false, $true, 1, $false, "547×3=1641", 2)
$d.Content.Find.Execute("775×4=3100", $true, $false, $false, $false, $false, $true, 1, $false, "194×7=1358", 2)
$d.Content.Find.Execute("838×7=5866", $true, $false, $false, $false, $false, $true, 1, $false, "329×2=658", 2)
$d.Content.Find.Execute("549×3=1647", $true, $false, $false, $false, $false, $true, 1, $false, "615×2=1230", 2)
$d.Content.Find.Execute("678×6=4068", $true, $false, $false, $false, $false, $true, 1, $false, "430×7=3010", 2)
$d.Content.Find.Execute("431×6=2586", $true, $false, $false, $false, $false, $true, 1, $false, "808×5=4040", 2)
$d.Content.Find.Execute("903×8=7224", $true, $false, $false, $false, $false, $true, 1, $false, "852×3=2556", 2)
$d.Content.Find.Execute("569×7=3983", $true, $false, $false, $false, $false, $true, 1, $false, "979×2=1958", 2)
$d.Content.Find.Execute("588×2=1176", $true, $false, $false, $false, $false, $true, 1, $false, "671×3=2013", 2)
